$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder "Periodo Mora" column (E16:E22) from descending (2409..2403) to
# ascending (2403..2409), and bring the "Valor Mora" (F) figure that
# travels with period 2409 along with it. All "Salario Basico" (G) values
# are updated to the new base salary of 1,300,000.

$periodos = @("2403", "2404", "2405", "2406", "2407", "2408", "2409")
$valorMora = @(52000, 52000, 52000, 52000, 52000, 52000, 32933)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = 1300000
}
